$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 9988
$ws.Range("D2").Value = 9988

$ws.Range("C3").Value = 9988
$ws.Range("D3").Value = 9988

$ws.Range("C4").Value = 9988
$ws.Range("D4").Value = 9985

$ws.Range("C5").Value = 9988
$ws.Range("D5").Value = 9988

$ws.Range("D6").Value = 9965

$ws.Range("C7").Value = 9988
$ws.Range("D7").Value = 9988

$ws.Range("C8").Value = 9988

$ws.Range("C9").Value = 5852
$ws.Range("D9").Value = 6969

$ws.Range("C10").Value = 5852
$ws.Range("D10").Value = 4437

$wb.Save()
